$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Propellant Combinations")
$ws.Activate()

# Remove the "Oxidizer CEA" and "Fuel CEA" columns (D:E) -- the CEA-specific
# propellant names are no longer needed; this shifts O:F start/stop/step left
# from F:H to D:F.
$ws.Range("D1:E1").EntireColumn.Delete()

# Leave the selection where the edit happened
$ws.Range("E11").Select()
